$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. ngs sheet: update existing "File transfer times" numbers
# ---------------------------------------------------------------------------
$ngs = $wb.Worksheets.Item("ngs")

$ngs.Range("K24").Value = 0.6875
$ngs.Range("K25").Value = 112
$ngs.Range("G26").Value = 26.91
$ngs.Range("K26").Value = 24

# ---------------------------------------------------------------------------
# 2. Add the new "file transfer" worksheet after the last existing sheet
# ---------------------------------------------------------------------------
$namd = $wb.Worksheets.Item("namd")
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $namd)
$ws3.Name = "file transfer"

$ws3.Range("C1").Value = " MB/sec avg"
$ws3.Range("E1").Value = " MB/sec avg"

$ws3.Range("C2").Value = 0.0833333333333333
$ws3.Range("E2").Value = 0.375
$ws3.Range("G2").Value = 0.6875
$ws3.Range("I2").Value = 0.5625
$ws3.Range("K2").Value = "old Speeds"

$ws3.Range("A3").Value = "Cyder to QB"
$ws3.Range("C3").Value = 113
$ws3.Range("E3").Value = 113.17
$ws3.Range("G3").Value = 112
$ws3.Range("I3").Value = 110
$ws3.Range("K3").Value = 70

$ws3.Range("A4").Value = "Cyder to Ranger"
$ws3.Range("C4").Value = 26.91
$ws3.Range("E4").Value = 23.6
$ws3.Range("G4").Value = 24.3
$ws3.Range("I4").Value = 24.2
$ws3.Range("K4").Value = 6

# ---------------------------------------------------------------------------
# 3. View state - selections on each sheet, and which sheet / cell is active
# ---------------------------------------------------------------------------
$namd.Select()
$namd.Range("F11").Select()

$ngs.Select()
$ngs.Range("I31").Select()

$ws3.Select()
$ws3.Range("H26").Select()
